$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "Golang Developer"
$ws.Range("B45").Value = "https://www.dice.com/job-detail/9e6797b8-138a-4d89-a76a-5dcefec2dccf"
$ws.Range("C45").Value = "Remote"
$ws.Range("D45").Value = "Contract"
$ws.Range("E45").Value = "Depends on Experience"
$ws.Range("F45").Value = "Source Mantra Inc"
